# Update the treeStats summary-statistics table with refreshed values
# (author re-ran the analysis and pasted in new results; a couple of
# CI bounds, point estimates, and one row label changed).
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("treeStats")

$ws.Range("B2").Value = "50.3 (25.1-26.9)"
$ws.Range("C2").Value = "48.1 (24.395-26.2)"
$ws.Range("D2").Value = "32.3 (18.4-20.5)"
$ws.Range("E2").Value = "36.7 (18.8-20.5)"
$ws.Range("B3").Value = "31.1 (14.1-16.2)"
$ws.Range("C3").Value = "28.8 (13.7-15.7)"
$ws.Range("D3").Value = "19.8 (10.7-12.8)"
$ws.Range("E3").Value = "33.8 (10.9-13)"
$ws.Range("B4").Value = "748 (880-917)"
$ws.Range("C4").Value = "633 (731-763)"
$ws.Range("D4").Value = "128 (159-175)"
$ws.Range("E4").Value = "155 (173.95-189)"
$ws.Range("B5").Value = "65914 (15212.85-22270.1)"
$ws.Range("C5").Value = "49519 (11861.55-18184.8)"
$ws.Range("D5").Value = "6990 (2005.05-3284.95)"
$ws.Range("E5").Value = "15605 (2190.95-3572.25)"
$ws.Range("B6").Value = "34 (79.9-146)"
$ws.Range("C6").Value = "26 (70-126)"
$ws.Range("E6").Value = "10 (22-44)"
$ws.Range("B7").Value = "1867 (1624.95-1687)"
$ws.Range("C7").Value = "1543 (1351-1401)"
$ws.Range("D7").Value = "355 (293-318)"
$ws.Range("E7").Value = "365 (318-346)"
$ws.Range("A8").Value = "Log roooted quartet index"
$ws.Range("C8").Value = "27.6 (27.3-27.9)"
$ws.Range("B9").Value = "83874 (38117.95-43968.05)"
$ws.Range("C9").Value = "64585 (30679.35-35479.55)"
$ws.Range("D9").Value = "9902 (5329.95-6376.2)"
$ws.Range("E9").Value = "18333 (5926.95-7001.8)"
$ws.Range("B10").Value = "0.692 (0.602-0.625)"
$ws.Range("C10").Value = "0.688 (0.60195-0.625)"
$ws.Range("D10").Value = "0.713 (0.584-0.6431)"
$ws.Range("E10").Value = "0.673 (0.5889-0.638)"
$ws.Range("B11").Value = "16.9 (15.4-16.2)"
$ws.Range("C11").Value = "16.3 (15-15.9)"
$ws.Range("D11").Value = "13 (12-12.8)"

# The author ended the session with treeStats active and B2:E11 selected,
# which is also reflected by the workbook-level activeTab bookkeeping.
$ws.Activate()
$ws.Range("B2:E11").Select()
